# "Sort stationing + size"
# Add a new sheet (placed first, named "Sheet1") containing the pull-sheet
# data sorted by stationing + cable size; rename the original "Sheet1" to
# "Sheet2.1" and move it to the second tab position.

$wb = $excel.ActiveWorkbook

# Worksheets.Add() inserts the new sheet before the active sheet, which is
# exactly the tab ordering we need (new sheet first, old sheet second).
$newSheet = $wb.Worksheets.Add()

# Grab a handle to the existing sheet after adding the new one.
$oldSheet = $wb.Worksheets.Item("Sheet1")

# Rename dance: free up the "Sheet1" name from the old sheet before giving
# it to the new one.
$newSheet.Name = "Sheet1New"
$oldSheet.Name = "Sheet2.1"
$newSheet.Name = "Sheet1"

$ws = $newSheet

$data = New-Object 'object[,]' 16,9

$data[0,0] = "Pull #"
$data[0,1] = "Local / Express"
$data[0,2] = "From "
$data[0,3] = "To"
$data[0,4] = "SK-#"
$data[0,5] = "Cable Size"
$data[0,6] = "Cable Type"
$data[0,7] = "Len. Before Mess."
$data[0,8] = "Len. After Mess."

$data[1,0] = 1;  $data[1,1] = "LOCAL"; $data[1,2] = "543+00"; $data[1,3] = "554+90"; $data[1,4] = 1; $data[1,5] = "7C#14";  $data[1,6] = "PK"; $data[1,7] = 0; $data[1,8] = 0
$data[2,0] = 2;  $data[2,1] = "LOCAL"; $data[2,2] = "543+00"; $data[2,3] = "554+90"; $data[2,4] = 1; $data[2,5] = "7C#14";  $data[2,6] = "PK"; $data[2,7] = 0; $data[2,8] = 0
$data[3,0] = 3;  $data[3,1] = "LOCAL"; $data[3,2] = "543+00"; $data[3,3] = "554+90"; $data[3,4] = 1; $data[3,5] = "7C#12";  $data[3,6] = "PK"; $data[3,7] = 0; $data[3,8] = 0
$data[4,0] = 4;  $data[4,1] = "LOCAL"; $data[4,2] = "543+00"; $data[4,3] = "553+00"; $data[4,4] = 2; $data[4,5] = "3C#6";   $data[4,6] = "PK"; $data[4,7] = 0; $data[4,8] = 0
$data[5,0] = 5;  $data[5,1] = "LOCAL"; $data[5,2] = "543+00"; $data[5,3] = "553+00"; $data[5,4] = 2; $data[5,5] = "12C#14"; $data[5,6] = "PK"; $data[5,7] = 0; $data[5,8] = 0
$data[6,0] = 6;  $data[6,1] = "LOCAL"; $data[6,2] = "543+00"; $data[6,3] = "553+00"; $data[6,4] = 2; $data[6,5] = "12C#14"; $data[6,6] = "PK"; $data[6,7] = 0; $data[6,8] = 0
$data[7,0] = 7;  $data[7,1] = "LOCAL"; $data[7,2] = "543+00"; $data[7,3] = "553+00"; $data[7,4] = 2; $data[7,5] = "12C#14"; $data[7,6] = "PK"; $data[7,7] = 0; $data[7,8] = 0
$data[8,0] = 8;  $data[8,1] = "LOCAL"; $data[8,2] = "543+00"; $data[8,3] = "553+00"; $data[8,4] = 2; $data[8,5] = "12C#14"; $data[8,6] = "PK"; $data[8,7] = 0; $data[8,8] = 0
$data[9,0] = 9;  $data[9,1] = "LOCAL"; $data[9,2] = "543+00"; $data[9,3] = "553+00"; $data[9,4] = 2; $data[9,5] = "12C#14"; $data[9,6] = "PK"; $data[9,7] = 0; $data[9,8] = 0
$data[10,0] = 10; $data[10,1] = "LOCAL"; $data[10,2] = "543+00"; $data[10,3] = "553+00"; $data[10,4] = 2; $data[10,5] = "7C#12";  $data[10,6] = "PK"; $data[10,7] = 0; $data[10,8] = 0
$data[11,0] = 11; $data[11,1] = "LOCAL"; $data[11,2] = "543+00"; $data[11,3] = "553+00"; $data[11,4] = 2; $data[11,5] = "12C#14"; $data[11,6] = "PK"; $data[11,7] = 0; $data[11,8] = 0
$data[12,0] = 12; $data[12,1] = "LOCAL"; $data[12,2] = "543+00"; $data[12,3] = "553+00"; $data[12,4] = 2; $data[12,5] = "12C#14"; $data[12,6] = "PK"; $data[12,7] = 0; $data[12,8] = 0
$data[13,0] = 13; $data[13,1] = "LOCAL"; $data[13,2] = "543+00"; $data[13,3] = "553+00"; $data[13,4] = 2; $data[13,5] = "12C#14"; $data[13,6] = "PK"; $data[13,7] = 0; $data[13,8] = 0
$data[14,0] = 14; $data[14,1] = "LOCAL"; $data[14,2] = "543+00"; $data[14,3] = "553+00"; $data[14,4] = 2; $data[14,5] = "12C#14"; $data[14,6] = "PK"; $data[14,7] = 0; $data[14,8] = 0
$data[15,0] = 15; $data[15,1] = "LOCAL"; $data[15,2] = "543+00"; $data[15,3] = "553+00"; $data[15,4] = 2; $data[15,5] = "7C#12";  $data[15,6] = "PK"; $data[15,7] = 0; $data[15,8] = 0

$ws.Range("A1:I16").Value = $data

# Match the author's final selection on the new sheet.
$ws.Range("D17").Select()
